$d = $word.ActiveDocument

# 1 & 8: Title heading and bold repeated title text (same replacement, applies globally)
$d.Content.Find.Execute("Play Halloween Jack Slot for Free | Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Halloween Jack - Free Slot Game Review", 2)

# 2: "Unique Wild symbol that awards free spins" -> "Detailed graphics and atmospheric visuals"
$d.Content.Find.Execute("Unique Wild symbol that awards free spins", $true, $false, $false, $false, $false, $true, 1, $false, "Detailed graphics and atmospheric visuals", 2)

# 3: "Great graphics with attention to detail" -> "Unique Wild symbol with free spins"
$d.Content.Find.Execute("Great graphics with attention to detail", $true, $false, $false, $false, $false, $true, 1, $false, "Unique Wild symbol with free spins", 2)

# 4: "Easy to get free spins" -> "Easy to trigger free spins compared to other slots"
$d.Content.Find.Execute("Easy to get free spins", $true, $false, $false, $false, $false, $true, 1, $false, "Easy to trigger free spins compared to other slots", 2)

# 5: "Atmospheric theme inspired by the horror genre" -> "Engaging theme inspired by horror and Halloween"
$d.Content.Find.Execute("Atmospheric theme inspired by the horror genre", $true, $false, $false, $false, $false, $true, 1, $false, "Engaging theme inspired by horror and Halloween", 2)

# 6: "Limited bonus features" -> "Limited bonus features beyond free spins"
$d.Content.Find.Execute("Limited bonus features", $true, $false, $false, $false, $false, $true, 1, $false, "Limited bonus features beyond free spins", 2)

# 7: "Low maximum bet" -> "May not appeal to players who are not fans of horror genre"
$d.Content.Find.Execute("Low maximum bet", $true, $false, $false, $false, $false, $true, 1, $false, "May not appeal to players who are not fans of horror genre", 2)

# 9: Meta description italic text
$d.Content.Find.Execute("Read our review of Halloween Jack slot game. Play for free with unique Wild symbol, great graphics, and easy free spins. Inspired by the horror genre.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Halloween Jack slot game and play it for free. Experience horror-themed slots with unique wild symbols and easy free spins.", 2)
